$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1: new "time_taken" column header, styled like the other headers (bold/centered/bordered)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Data cells F2:F48: per-row "time_taken" metadata timestamps
$timeTaken = @(
    "2021-10-05 13:39:37.775339",
    "2021-10-05 13:39:37.775352",
    "2021-10-05 13:39:37.775357",
    "2021-10-05 13:39:37.775360",
    "2021-10-05 13:39:37.775363",
    "2021-10-05 13:39:37.775367",
    "2021-10-05 13:39:37.775370",
    "2021-10-05 13:39:37.775373",
    "2021-10-05 13:39:37.775376",
    "2021-10-05 13:39:37.775379",
    "2021-10-05 13:39:37.775382",
    "2021-10-05 13:39:37.775385",
    "2021-10-05 13:39:37.775392",
    "2021-10-05 13:39:37.775395",
    "2021-10-05 13:39:37.775398",
    "2021-10-05 13:39:37.775401",
    "2021-10-05 13:39:37.775424",
    "2021-10-05 13:39:37.775428",
    "2021-10-05 13:39:37.775431",
    "2021-10-05 13:39:37.775434",
    "2021-10-05 13:39:37.775436",
    "2021-10-05 13:39:37.775439",
    "2021-10-05 13:39:37.775441",
    "2021-10-05 13:39:37.775444",
    "2021-10-05 13:39:37.775447",
    "2021-10-05 13:39:37.775450",
    "2021-10-05 13:39:37.775452",
    "2021-10-05 13:39:37.775454",
    "2021-10-05 13:39:37.775457",
    "2021-10-05 13:39:37.775459",
    "2021-10-05 13:39:37.775462",
    "2021-10-05 13:39:37.775464",
    "2021-10-05 13:39:37.775467",
    "2021-10-05 13:39:37.775470",
    "2021-10-05 13:39:37.775473",
    "2021-10-05 13:39:37.775475",
    "2021-10-05 13:39:37.775478",
    "2021-10-05 13:39:37.775480",
    "2021-10-05 13:39:37.775483",
    "2021-10-05 13:39:37.775485",
    "2021-10-05 13:39:37.775488",
    "2021-10-05 13:39:37.775491",
    "2021-10-05 13:39:37.775494",
    "2021-10-05 13:39:37.775496",
    "2021-10-05 13:39:37.775499",
    "2021-10-05 13:39:37.775501",
    "2021-10-05 13:39:37.775504"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}
